# Update res_bus/vm_pu.xlsx values for the 380 kV case (rows 2-25).
# Column B is the slack-bus voltage setpoint (1.05 -> 1.02 p.u.);
# columns C-F and I-N are the recomputed per-bus voltage magnitudes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.027644874290315
$ws.Cells.Item(2, 4).Value = 1.035205882930603
$ws.Cells.Item(2, 5).Value = 1.03635974369666
$ws.Cells.Item(2, 6).Value = 1.044138801283459
$ws.Cells.Item(2, 9).Value = 1.034443573409688
$ws.Cells.Item(2, 10).Value = 1.032801598442629
$ws.Cells.Item(2, 11).Value = 1.038003153577268
$ws.Cells.Item(2, 12).Value = 1.039153709507249
$ws.Cells.Item(2, 13).Value = 1.046910693532092
$ws.Cells.Item(2, 14).Value = 1.014945071428707
# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.028481254863815
$ws.Cells.Item(3, 4).Value = 1.035685638156872
$ws.Cells.Item(3, 5).Value = 1.037137318184863
$ws.Cells.Item(3, 6).Value = 1.045077543415435
$ws.Cells.Item(3, 9).Value = 1.034572099891409
$ws.Cells.Item(3, 10).Value = 1.033278832566549
$ws.Cells.Item(3, 11).Value = 1.038293132674698
$ws.Cells.Item(3, 12).Value = 1.039740954254256
$ws.Cells.Item(3, 13).Value = 1.04766027984922
$ws.Cells.Item(3, 14).Value = 1.015104732275649
# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.029023068487988
$ws.Cells.Item(4, 4).Value = 1.035996276719077
$ws.Cells.Item(4, 5).Value = 1.03764145916311
$ws.Cells.Item(4, 6).Value = 1.04568632754884
$ws.Cells.Item(4, 9).Value = 1.034654235789662
$ws.Cells.Item(4, 10).Value = 1.033587603756805
$ws.Cells.Item(4, 11).Value = 1.038480250832825
$ws.Cells.Item(4, 12).Value = 1.040121279854186
$ws.Cells.Item(4, 13).Value = 1.048146040815626
$ws.Cells.Item(4, 14).Value = 1.015207988811583
# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.029250993684819
$ws.Cells.Item(5, 4).Value = 1.036126916247814
$ws.Cells.Item(5, 5).Value = 1.037853637282363
$ws.Cells.Item(5, 6).Value = 1.045942583043806
$ws.Cells.Item(5, 9).Value = 1.034688518566463
$ws.Cells.Item(5, 10).Value = 1.033717402497613
$ws.Cells.Item(5, 11).Value = 1.038558790029763
$ws.Cells.Item(5, 12).Value = 1.040281248518145
$ws.Cells.Item(5, 13).Value = 1.048350427559363
$ws.Cells.Item(5, 14).Value = 1.01525138431343
# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.029289271888686
$ws.Cells.Item(6, 4).Value = 1.036148853902914
$ws.Cells.Item(6, 5).Value = 1.037889276784193
$ws.Cells.Item(6, 6).Value = 1.045985628337003
$ws.Cells.Item(6, 9).Value = 1.034694260281017
$ws.Cells.Item(6, 10).Value = 1.033739195734245
$ws.Cells.Item(6, 11).Value = 1.038571969721318
$ws.Cells.Item(6, 12).Value = 1.040308112599325
$ws.Cells.Item(6, 13).Value = 1.048384755082423
$ws.Cells.Item(6, 14).Value = 1.01525866980176
# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.029026113460936
$ws.Cells.Item(7, 4).Value = 1.035998022148356
$ws.Cells.Item(7, 5).Value = 1.037644293367471
$ws.Cells.Item(7, 6).Value = 1.045689750381775
$ws.Cells.Item(7, 9).Value = 1.034654694849867
$ws.Cells.Item(7, 10).Value = 1.033589338169492
$ws.Cells.Item(7, 11).Value = 1.038481300771354
$ws.Cells.Item(7, 12).Value = 1.040123417051533
$ws.Cells.Item(7, 13).Value = 1.048148771164148
$ws.Cells.Item(7, 14).Value = 1.015208568718192
# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.027927403822001
$ws.Cells.Item(8, 4).Value = 1.035367974860465
$ws.Cells.Item(8, 5).Value = 1.036622320886178
$ws.Cells.Item(8, 6).Value = 1.044455771899304
$ws.Cells.Item(8, 9).Value = 1.034487222225548
$ws.Cells.Item(8, 10).Value = 1.032962887688371
$ws.Cells.Item(8, 11).Value = 1.038101259763078
$ws.Cells.Item(8, 12).Value = 1.039352100605902
$ws.Cells.Item(8, 13).Value = 1.047163867952588
$ws.Cells.Item(8, 14).Value = 1.014999040617419
# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02599615016252
$ws.Cells.Item(9, 4).Value = 1.034259406641556
$ws.Cells.Item(9, 5).Value = 1.034829197491867
$ws.Cells.Item(9, 6).Value = 1.042291802750899
$ws.Cells.Item(9, 9).Value = 1.034184262983883
$ws.Cells.Item(9, 10).Value = 1.031858815953945
$ws.Cells.Item(9, 11).Value = 1.037427672932877
$ws.Cells.Item(9, 12).Value = 1.037995599314076
$ws.Cells.Item(9, 13).Value = 1.045433991544671
$ws.Cells.Item(9, 14).Value = 1.014629426160579
# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.024711976206868
$ws.Cells.Item(10, 4).Value = 1.033521591576209
$ws.Cells.Item(10, 5).Value = 1.033639078050204
$ws.Cells.Item(10, 6).Value = 1.040856295377704
$ws.Cells.Item(10, 9).Value = 1.033977055238756
$ws.Cells.Item(10, 10).Value = 1.031122714541544
$ws.Cells.Item(10, 11).Value = 1.036976068645545
$ws.Cells.Item(10, 12).Value = 1.037093133182446
$ws.Cells.Item(10, 13).Value = 1.04428463224098
$ws.Cells.Item(10, 14).Value = 1.01438277486907
# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024156723586993
$ws.Cells.Item(11, 4).Value = 1.033202426224389
$ws.Cells.Item(11, 5).Value = 1.033125021297423
$ws.Cells.Item(11, 6).Value = 1.04023642059
$ws.Cells.Item(11, 9).Value = 1.033886099573641
$ws.Cells.Item(11, 10).Value = 1.030803976841915
$ws.Cells.Item(11, 11).Value = 1.036779931929405
$ws.Cells.Item(11, 12).Value = 1.036702815701724
$ws.Cells.Item(11, 13).Value = 1.043787889753148
$ws.Cells.Item(11, 14).Value = 1.014275920453883
# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.023950600262207
$ws.Cells.Item(12, 4).Value = 1.03308392312026
$ws.Cells.Item(12, 5).Value = 1.032934270733959
$ws.Cells.Item(12, 6).Value = 1.040006430130626
$ws.Cells.Item(12, 9).Value = 1.033852130025696
$ws.Cells.Item(12, 10).Value = 1.030685584535623
$ws.Cells.Item(12, 11).Value = 1.036706990882854
$ws.Cells.Item(12, 12).Value = 1.036557904353841
$ws.Cells.Item(12, 13).Value = 1.043603519941482
$ws.Cells.Item(12, 14).Value = 1.014236222489321
# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.023994808859357
$ws.Cells.Item(13, 4).Value = 1.033109340183177
$ws.Cells.Item(13, 5).Value = 1.032975178615034
$ws.Cells.Item(13, 6).Value = 1.040055752131747
$ws.Cells.Item(13, 9).Value = 1.033859424949202
$ws.Cells.Item(13, 10).Value = 1.030710980019538
$ws.Cells.Item(13, 11).Value = 1.036722640906199
$ws.Cells.Item(13, 12).Value = 1.036588985142499
$ws.Cells.Item(13, 13).Value = 1.043643061419508
$ws.Cells.Item(13, 14).Value = 1.014244738168571
# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024139682866404
$ws.Cells.Item(14, 4).Value = 1.033192629707136
$ws.Cells.Item(14, 5).Value = 1.03310924984887
$ws.Cells.Item(14, 6).Value = 1.040217404220226
$ws.Cells.Item(14, 9).Value = 1.033883295401029
$ws.Cells.Item(14, 10).Value = 1.030794190468192
$ws.Cells.Item(14, 11).Value = 1.036773904369431
$ws.Cells.Item(14, 12).Value = 1.03669083584965
$ws.Cells.Item(14, 13).Value = 1.043772646771537
$ws.Cells.Item(14, 14).Value = 1.014272639157663
# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024228960765552
$ws.Cells.Item(15, 4).Value = 1.033243953700925
$ws.Cells.Item(15, 5).Value = 1.033191881188831
$ws.Cells.Item(15, 6).Value = 1.040317037733547
$ws.Cells.Item(15, 9).Value = 1.033897978339588
$ws.Cells.Item(15, 10).Value = 1.030845459349421
$ws.Cells.Item(15, 11).Value = 1.036805477975293
$ws.Cells.Item(15, 12).Value = 1.036753598741527
$ws.Cells.Item(15, 13).Value = 1.04385250750949
$ws.Cells.Item(15, 14).Value = 1.014289828900633
# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.024748843549023
$ws.Cells.Item(16, 4).Value = 1.033542780306528
$ws.Cells.Item(16, 5).Value = 1.033673221272332
$ws.Cells.Item(16, 6).Value = 1.040897470595948
$ws.Cells.Item(16, 9).Value = 1.033983065744012
$ws.Cells.Item(16, 10).Value = 1.031143868196333
$ws.Cells.Item(16, 11).Value = 1.036989073298829
$ws.Cells.Item(16, 12).Value = 1.037119047000539
$ws.Cells.Item(16, 13).Value = 1.04431761929539
$ws.Cells.Item(16, 14).Value = 1.014389865371322
# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.025075168412073
$ws.Cells.Item(17, 4).Value = 1.033730311691791
$ws.Cells.Item(17, 5).Value = 1.033975495213552
$ws.Cells.Item(17, 6).Value = 1.041262019730084
$ws.Cells.Item(17, 9).Value = 1.034036109140066
$ws.Cells.Item(17, 10).Value = 1.031331052784801
$ws.Cells.Item(17, 11).Value = 1.037104081046347
$ws.Cells.Item(17, 12).Value = 1.037348406185006
$ws.Cells.Item(17, 13).Value = 1.044609623795185
$ws.Cells.Item(17, 14).Value = 1.014452601795542
# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.025265585574328
$ws.Cells.Item(18, 4).Value = 1.033839725682905
$ws.Cells.Item(18, 5).Value = 1.034151929178854
$ws.Cells.Item(18, 6).Value = 1.041474819880178
$ws.Cells.Item(18, 9).Value = 1.034066929406624
$ws.Cells.Item(18, 10).Value = 1.031440234217573
$ws.Cells.Item(18, 11).Value = 1.037171106190676
$ws.Cells.Item(18, 12).Value = 1.037482231450988
$ws.Cells.Item(18, 13).Value = 1.044780035444872
$ws.Cells.Item(18, 14).Value = 1.01448918976771
# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.025330525985451
$ws.Cells.Item(19, 4).Value = 1.033877038071221
$ws.Cells.Item(19, 5).Value = 1.034212109357783
$ws.Cells.Item(19, 6).Value = 1.041547407130988
$ws.Cells.Item(19, 9).Value = 1.03407741810603
$ws.Cells.Item(19, 10).Value = 1.031477462198974
$ws.Cells.Item(19, 11).Value = 1.037193950335069
$ws.Cells.Item(19, 12).Value = 1.03752786981811
$ws.Cells.Item(19, 13).Value = 1.044838156692739
$ws.Cells.Item(19, 14).Value = 1.014501664438046
# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.025040148826494
$ws.Cells.Item(20, 4).Value = 1.033710188220209
$ws.Cells.Item(20, 5).Value = 1.033943051377189
$ws.Cells.Item(20, 6).Value = 1.041222889998054
$ws.Cells.Item(20, 9).Value = 1.03403043039418
$ws.Cells.Item(20, 10).Value = 1.031310969658908
$ws.Cells.Item(20, 11).Value = 1.037091747682046
$ws.Cells.Item(20, 12).Value = 1.037323793560816
$ws.Cells.Item(20, 13).Value = 1.044578285117475
$ws.Cells.Item(20, 14).Value = 1.014445871297551
# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024097017695836
$ws.Cells.Item(21, 4).Value = 1.033168101646114
$ws.Cells.Item(21, 5).Value = 1.033069763870054
$ws.Cells.Item(21, 6).Value = 1.040169794554957
$ws.Cells.Item(21, 9).Value = 1.033876271236914
$ws.Cells.Item(21, 10).Value = 1.030769687025874
$ws.Cells.Item(21, 11).Value = 1.036758810951683
$ws.Cells.Item(21, 12).Value = 1.036660841413027
$ws.Cells.Item(21, 13).Value = 1.043734483169658
$ws.Cells.Item(21, 14).Value = 1.014264423214562
# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023504741624908
$ws.Cells.Item(22, 4).Value = 1.032827555483601
$ws.Cells.Item(22, 5).Value = 1.032521810888144
$ws.Cells.Item(22, 6).Value = 1.03950916894865
$ws.Cells.Item(22, 9).Value = 1.033778277787792
$ws.Cells.Item(22, 10).Value = 1.030429367923238
$ws.Cells.Item(22, 11).Value = 1.036548976939518
$ws.Cells.Item(22, 12).Value = 1.036244422443507
$ws.Cells.Item(22, 13).Value = 1.043204776446514
$ws.Cells.Item(22, 14).Value = 1.014150296566863
# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.023818650755763
$ws.Cells.Item(23, 4).Value = 1.033008057740679
$ws.Cells.Item(23, 5).Value = 1.032812184608989
$ws.Cells.Item(23, 6).Value = 1.039859236546563
$ws.Cells.Item(23, 9).Value = 1.033830326908401
$ws.Cells.Item(23, 10).Value = 1.03060977647345
$ws.Cells.Item(23, 11).Value = 1.036660261188219
$ws.Cells.Item(23, 12).Value = 1.036465135131186
$ws.Cells.Item(23, 13).Value = 1.043485505265792
$ws.Cells.Item(23, 14).Value = 1.014210801190526
# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.025055972430873
$ws.Cells.Item(24, 4).Value = 1.033719281058061
$ws.Cells.Item(24, 5).Value = 1.033957710971447
$ws.Cells.Item(24, 6).Value = 1.041240570530632
$ws.Cells.Item(24, 9).Value = 1.034032996742941
$ws.Cells.Item(24, 10).Value = 1.031320044359705
$ws.Cells.Item(24, 11).Value = 1.037097320774572
$ws.Cells.Item(24, 12).Value = 1.03733491481067
$ws.Cells.Item(24, 13).Value = 1.044592445437954
$ws.Cells.Item(24, 14).Value = 1.014448912535766
# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.02649484563872
$ws.Cells.Item(25, 4).Value = 1.034545789791311
$ws.Cells.Item(25, 5).Value = 1.035291836991418
$ws.Cells.Item(25, 6).Value = 1.042849990237373
$ws.Cells.Item(25, 9).Value = 1.034263511102194
$ws.Cells.Item(25, 10).Value = 1.032144259640188
$ws.Cells.Item(25, 11).Value = 1.038003153577268
$ws.Cells.Item(25, 12).Value = 1.038345964143762
$ws.Cells.Item(25, 13).Value = 1.045880526997339
$ws.Cells.Item(25, 14).Value = 1.014725024906708
